# Apply numeric corrections to the Masamune_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each sheet holds a "Leve" profit table; this patches currentAveragePrice / LevePrice / LeveProfit
# columns (H, I, J, K, L, M, N) for the affected rows.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 49667
$ws.Range("J105").Value = 49667
$ws.Range("L105").Value = 49667
$ws.Range("N105").Value = -56655
$ws.Range("H109").Value = 32113.5
$ws.Range("J109").Value = 32113.5
$ws.Range("L109").Value = 32113.5
$ws.Range("N109").Value = -34887.5
$ws.Range("H110").Value = 25550.5
$ws.Range("J110").Value = 25550.5
$ws.Range("L110").Value = 25550.5
$ws.Range("N110").Value = -33730.5
$ws.Range("H117").Value = 48538
$ws.Range("J117").Value = 48538
$ws.Range("L117").Value = 48538
$ws.Range("N117").Value = -57716

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 29240
$ws.Range("J37").Value = 29240
$ws.Range("L37").Value = 29240
$ws.Range("N37").Value = -29786
$ws.Range("H75").Value = 21578.5
$ws.Range("J75").Value = 33000
$ws.Range("L75").Value = 33000
$ws.Range("N75").Value = -34748
$ws.Range("H78").Value = 21578.5
$ws.Range("J78").Value = 33000
$ws.Range("L78").Value = 99000
$ws.Range("N78").Value = -107736
$ws.Range("H80").Value = 59106
$ws.Range("J80").Value = 59106
$ws.Range("L80").Value = 59106
$ws.Range("N80").Value = -61102
$ws.Range("H83").Value = 59106
$ws.Range("J83").Value = 59106
$ws.Range("L83").Value = 177318
$ws.Range("N83").Value = -187302
$ws.Range("H120").Value = 43296
$ws.Range("J120").Value = 43296
$ws.Range("L120").Value = 43296
$ws.Range("N120").Value = -52972

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J5").Value = 500
$ws.Range("L5").Value = 500
$ws.Range("N5").Value = -726
$ws.Range("H10").Value = 700
$ws.Range("I10").Value = 300
$ws.Range("J10").Value = 900
$ws.Range("K10").Value = 300
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = -160
$ws.Range("N10").Value = -1180
$ws.Range("H76").Value = 24000
$ws.Range("J76").Value = 24000
$ws.Range("L76").Value = 24000
$ws.Range("N76").Value = -24630
$ws.Range("H79").Value = 24000
$ws.Range("J79").Value = 24000
$ws.Range("L79").Value = 24000
$ws.Range("N79").Value = -26184
$ws.Range("H117").Value = 48933.5
$ws.Range("J117").Value = 48933.5
$ws.Range("L117").Value = 48933.5
$ws.Range("N117").Value = -58111.5
$ws.Range("H130").Value = 40552.855
$ws.Range("J130").Value = 40552.855
$ws.Range("L130").Value = 40552.855
$ws.Range("N130").Value = -50592.855

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 420.58334
$ws.Range("I5").Value = 316.66666
$ws.Range("J5").Value = 455.22223
$ws.Range("K5").Value = 316.66666
$ws.Range("L5").Value = 455.22223
$ws.Range("M5").Value = -204.66666
$ws.Range("N5").Value = -679.2222300000001
$ws.Range("H19").Value = 69
$ws.Range("I19").Value = 69
$ws.Range("K19").Value = 69
$ws.Range("M19").Value = 101
$ws.Range("H24").Value = 69
$ws.Range("I24").Value = 69
$ws.Range("K24").Value = 69
$ws.Range("M24").Value = 101
$ws.Range("H31").Value = 217252.72
$ws.Range("J31").Value = 255543.97
$ws.Range("L31").Value = 255543.97
$ws.Range("N31").Value = -256133.97
$ws.Range("H34").Value = 217252.72
$ws.Range("J34").Value = 255543.97
$ws.Range("L34").Value = 255543.97
$ws.Range("N34").Value = -255947.97

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 53300
$ws.Range("I3").Value = 76700
$ws.Range("J3").Value = 6500
$ws.Range("K3").Value = 76700
$ws.Range("L3").Value = 6500
$ws.Range("M3").Value = -76584
$ws.Range("N3").Value = -6732
$ws.Range("H5").Value = 15818.182
$ws.Range("J5").Value = 15900
$ws.Range("L5").Value = 15900
$ws.Range("N5").Value = -16124
$ws.Range("H6").Value = 20998
$ws.Range("I6").Value = 25000
$ws.Range("J6").Value = 19997.5
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 19997.5
$ws.Range("M6").Value = -24887
$ws.Range("N6").Value = -20223.5
$ws.Range("H9").Value = 2300
$ws.Range("I9").Value = 2300
$ws.Range("K9").Value = 2300
$ws.Range("M9").Value = -2130
$ws.Range("H16").Value = 20998
$ws.Range("I16").Value = 25000
$ws.Range("J16").Value = 19997.5
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 19997.5
$ws.Range("M16").Value = -24750
$ws.Range("N16").Value = -20497.5
$ws.Range("H110").Value = 32754
$ws.Range("J110").Value = 32754
$ws.Range("L110").Value = 32754
$ws.Range("N110").Value = -40934
$ws.Range("H122").Value = 2733.3333
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -11200

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 2394.8
$ws.Range("I9").Value = 235.42857
$ws.Range("J9").Value = 7433.3335
$ws.Range("K9").Value = 235.42857
$ws.Range("L9").Value = 7433.3335
$ws.Range("M9").Value = -11.42857000000001
$ws.Range("N9").Value = -7881.3335
$ws.Range("H104").Value = 16916.9
$ws.Range("J104").Value = 16916.9
$ws.Range("L104").Value = 16916.9
$ws.Range("N104").Value = -23904.9
$ws.Range("H105").Value = 33138
$ws.Range("J105").Value = 33138
$ws.Range("L105").Value = 33138
$ws.Range("N105").Value = -40126
$ws.Range("H106").Value = 32090.25
$ws.Range("J106").Value = 32090.25
$ws.Range("L106").Value = 32090.25
$ws.Range("N106").Value = -34614.25
$ws.Range("H110").Value = 41998.5
$ws.Range("J110").Value = 41998.5
$ws.Range("L110").Value = 41998.5
$ws.Range("N110").Value = -50178.5
$ws.Range("H111").Value = 38305
$ws.Range("J111").Value = 38305
$ws.Range("L111").Value = 38305
$ws.Range("N111").Value = -46485
$ws.Range("H112").Value = 26654.666
$ws.Range("J112").Value = 28985.6
$ws.Range("L112").Value = 28985.6
$ws.Range("N112").Value = -31939.6
$ws.Range("H115").Value = 30302
$ws.Range("J115").Value = 30302
$ws.Range("L115").Value = 30302
$ws.Range("N115").Value = -32652
$ws.Range("H121").Value = 27660
$ws.Range("J121").Value = 27660
$ws.Range("L121").Value = 27660
$ws.Range("N121").Value = -31154
$ws.Range("H124").Value = 44714.5
$ws.Range("J124").Value = 44714.5
$ws.Range("L124").Value = 44714.5
$ws.Range("N124").Value = -54534.5
$ws.Range("H125").Value = 49707
$ws.Range("J125").Value = 49707
$ws.Range("L125").Value = 49707
$ws.Range("N125").Value = -59547

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 93338
$ws.Range("I9").Value = 80000
$ws.Range("K9").Value = 80000
$ws.Range("M9").Value = -79860
$ws.Range("H119").Value = 46497
$ws.Range("J119").Value = 46497
$ws.Range("L119").Value = 46497
$ws.Range("N119").Value = -56173

